$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Give the new row 5 the same cell formatting as row 4 (e.g. the bordered/bold/centered
#    style used on column A) WITHOUT copying any values, so the shared-string table stays clean.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null

# 2) Clear the shared-string cells whose *content* changes, so the stale strings are dropped
#    from the workbook string table before we write the replacement text (keeps the string
#    table compact/ordered the same way Excel would rebuild it).
$ws.Range("C3:C4").Value = ""
$ws.Range("O2:T4").Value = ""

# 3) Re-populate the model-predictor list strings (column C) for rows 3-5 in order.
$ws.Range("C3").Value = "['N1ratio-ArgsPreds', 'latitude', 'longitude', 'Macro_class']"
$ws.Range("C4").Value = "['N1ratio-ArgsPreds', 'latitude', 'longitude', 'Macro_class', 'Fam_class']"
$ws.Range("C5").Value = "['N1ratio-ArgsPreds', 'latitude', 'longitude', 'Macro_class', 'Fam_class', 'Nlen_freq', 'Vlen_freq']"

# 4) Re-populate the coefficient/P-value dict strings (columns O-T) for rows 2-5, column by column.
$ws.Range("O2").Value = "{'const': 0.8769068789051704, 'N1ratio-ArgsPreds': -0.25892073248774317}"
$ws.Range("O3").Value = "{'const': 0.492442807410804, 'N1ratio-ArgsPreds': -0.2347084757011529, 'latitude': -0.00039709861635979916, 'longitude': -0.0005177668538437926, 'Macro_class': 0.13581942809818}"
$ws.Range("O4").Value = "{'const': 0.6390956805391714, 'N1ratio-ArgsPreds': -0.22578922616158928, 'latitude': -0.000525123248977932, 'longitude': -0.0010411331253083384, 'Macro_class': 0.11757684971827578, 'Fam_class': -0.002279536807291503}"
$ws.Range("O5").Value = "{'const': 0.07174266864899048, 'N1ratio-ArgsPreds': -0.22105968489772074, 'latitude': 0.0002667620768122055, 'longitude': -0.0013574861844785654, 'Macro_class': 0.10147993435936434, 'Fam_class': -0.0047730797442657015, 'Nlen_freq': 0.06963180770116814, 'Vlen_freq': 0.0385371271582626}"

$ws.Range("P2").Value = "{'const': 3.726577191319948e-15, 'N1ratio-ArgsPreds': 7.810452239194086e-10}"
$ws.Range("P3").Value = "{'const': 7.325298994937364e-05, 'N1ratio-ArgsPreds': 1.3185501625227715e-09, 'latitude': 0.787682088353225, 'longitude': 0.24372793217946584, 'Macro_class': 2.2147465213834944e-06}"
$ws.Range("P4").Value = "{'const': 0.0002704792054073951, 'N1ratio-ArgsPreds': 7.841109645115831e-09, 'latitude': 0.7218232433499598, 'longitude': 0.09398637531538649, 'Macro_class': 0.00024134675563298213, 'Fam_class': 0.22666500312154214}"
$ws.Range("P5").Value = "{'const': 0.7972316822536156, 'N1ratio-ArgsPreds': 1.2220126715507334e-08, 'latitude': 0.8567082391104337, 'longitude': 0.029702017900401882, 'Macro_class': 0.0013936046111223586, 'Fam_class': 0.025601998533392043, 'Nlen_freq': 0.23396876351533552, 'Vlen_freq': 0.3983691468942224}"

$ws.Range("Q2").Value = "{'N1ratio-ArgsPreds': -0.5531368928814658}"
$ws.Range("Q3").Value = "{'N1ratio-ArgsPreds': -0.5014118249044676, 'latitude': -0.022114931134383408, 'longitude': -0.10051783646596273, 'Macro_class': 0.39915299056925585}"
$ws.Range("Q4").Value = "{'N1ratio-ArgsPreds': -0.48235747599333034, 'latitude': -0.02924478708756917, 'longitude': -0.2021227285063198, 'Macro_class': 0.3455407804606315, 'Fam_class': -0.157703905320072}"
$ws.Range("Q5").Value = "{'N1ratio-ArgsPreds': -0.47225367420691, 'latitude': 0.014856322119797389, 'longitude': -0.26353864347096045, 'Macro_class': 0.2982343531371035, 'Fam_class': -0.3302132756387528, 'Nlen_freq': 0.14132340960025755, 'Vlen_freq': 0.09172628232179939}"

$ws.Range("R2").Value = "{'N1ratio-ArgsPreds': -0.5531368928814663}"
$ws.Range("R3").Value = "{'N1ratio-ArgsPreds': -0.5536010432239684, 'latitude': -0.026860080239888274, 'longitude': -0.11588859139307682, 'Macro_class': 0.4469540976552578}"
$ws.Range("R4").Value = "{'N1ratio-ArgsPreds': -0.5333790629645466, 'latitude': -0.035680055520583905, 'longitude': -0.16671544191102688, 'Macro_class': 0.35591767925655415, 'Fam_class': -0.12075793785607132}"
$ws.Range("R5").Value = "{'N1ratio-ArgsPreds': -0.5320484327264897, 'latitude': 0.018284980288520318, 'longitude': -0.21752858914027481, 'Macro_class': 0.3153651815706621, 'Fam_class': -0.22320178764535017, 'Nlen_freq': 0.12009834319405689, 'Vlen_freq': 0.08537275883833031}"

$ws.Range("S2").Value = "{'N1ratio-ArgsPreds': -0.5531368928814662}"
$ws.Range("S3").Value = "{'N1ratio-ArgsPreds': -0.49437227359470837, 'latitude': -0.019982626440729546, 'longitude': -0.08676914369418048, 'Macro_class': 0.37157244508480686}"
$ws.Range("S4").Value = "{'N1ratio-ArgsPreds': -0.46550866799422724, 'latitude': -0.026357294175478444, 'longitude': -0.12482321704588047, 'Macro_class': 0.2811647640366989, 'Fam_class': -0.08980576808409556}"
$ws.Range("S5").Value = "{'N1ratio-ArgsPreds': -0.44915676676172267, 'latitude': 0.01307225734629092, 'longitude': -0.159303749956017, 'Macro_class': 0.2375442705402404, 'Fam_class': -0.16367335530192872, 'Nlen_freq': 0.08647196942162136, 'Vlen_freq': 0.06124790955740818}"

$ws.Range("T2").Value = "{'N1ratio-ArgsPreds': 30.596042226656255}"
$ws.Range("T3").Value = "{'N1ratio-ArgsPreds': 24.440394489920116, 'latitude': 0.03993053594697436, 'longitude': 0.7528884297421341, 'Macro_class': 13.80660819463018}"
$ws.Range("T4").Value = "{'N1ratio-ArgsPreds': 21.66983199777597, 'latitude': 0.069470695625271, 'longitude': 1.5580835513682985, 'Macro_class': 7.905362453581256, 'Fam_class': 0.8065075981174357}"
$ws.Range("T5").Value = "{'N1ratio-ArgsPreds': 20.174180112784455, 'latitude': 0.017088391212765692, 'longitude': 2.537768475004919, 'Macro_class': 5.642728046649493, 'Fam_class': 2.67889672357914, 'Nlen_freq': 0.7477401495653819, 'Vlen_freq': 0.3751306425152452}"

# 5) Update the numeric statistics cells for the existing rows (2-4) to the refreshed values.
$ws.Range("D2").Value = 106
$ws.Range("E2").Value = 104
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = [double]"0.3059604222665624"
$ws.Range("H2").Value = [double]"45.84736222052121"
$ws.Range("I2").Value = [double]"7.810452239193906e-10"
$ws.Range("J2").Value = [double]"14.92839846445508"
$ws.Range("K2").Value = [double]"21.50943396226416"
$ws.Range("L2").Value = [double]"6.581035497809081"
$ws.Range("M2").Value = [double]"0.1435422929274527"
$ws.Range("N2").Value = [double]"0.2048517520215634"

$ws.Range("D3").Value = 106
$ws.Range("E3").Value = 101
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = [double]"0.4469338392920367"
$ws.Range("H3").Value = [double]"20.4045740706288"
$ws.Range("I3").Value = [double]"2.413267169335551e-12"
$ws.Range("J3").Value = [double]"11.89614006051091"
$ws.Range("K3").Value = [double]"21.50943396226416"
$ws.Range("L3").Value = [double]"2.403323475438311"
$ws.Range("M3").Value = [double]"0.1177835649555536"
$ws.Range("N3").Value = [double]"0.2048517520215634"
$ws.Range("U3").Value = [double]"0.1409734170254743"
$ws.Range("V3").Value = [double]"8.581441746105543"
$ws.Range("W3").Value = [double]"3.955623329731498e-05"

$ws.Range("D4").Value = 106
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = [double]"0.454998915273213"
$ws.Range("H4").Value = [double]"16.69717466714061"
$ws.Range("I4").Value = [double]"5.794027875875511e-12"
$ws.Range("J4").Value = [double]"11.72266484129316"
$ws.Range("K4").Value = [double]"21.50943396226416"
$ws.Range("L4").Value = [double]"1.9573538241942"
$ws.Range("M4").Value = [double]"0.1172266484129316"
$ws.Range("N4").Value = [double]"0.2048517520215634"
$ws.Range("U4").Value = [double]"0.00806507598117634"
$ws.Range("V4").Value = [double]"1.479827509924943"
$ws.Range("W4").Value = [double]"0.2266650031215442"

# 6) Fill in all values for the newly added row 5.
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 4
$ws.Range("D5").Value = 106
$ws.Range("E5").Value = 98
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = [double]"0.4890634649576791"
$ws.Range("H5").Value = [double]"13.40066336974782"
$ws.Range("I5").Value = [double]"4.989964903313639e-12"
$ws.Range("J5").Value = [double]"10.98995565940087"
$ws.Range("K5").Value = [double]"21.50943396226416"
$ws.Range("L5").Value = [double]"1.502782614694755"
$ws.Range("M5").Value = [double]"0.112142404687764"
$ws.Range("N5").Value = [double]"0.2048517520215634"
$ws.Range("U5").Value = [double]"0.03406454968446604"
$ws.Range("V5").Value = [double]"3.266869405611362"
$ws.Range("W5").Value = [double]"0.042316525175481"

